$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 115 is brand new; copy formatting for style-bearing columns A (id) and D (date) from row 114
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("D114").Copy()
$ws.Range("D115").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 112
$ws.Range("A112").Value = 110
$ws.Range("B112").Value = 7802944
$ws.Range("C112").Value = "Canada Premier League"
$ws.Range("D112").Value = 45436.83333333334
$ws.Range("E112").Value = "York United FC"
$ws.Range("F112").Value = "HFX Wanderers"
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 1
$ws.Range("I112").Value = "H"
$ws.Range("J112").Value = 1.909
$ws.Range("K112").Value = 3.4
$ws.Range("L112").Value = 3.4
$ws.Range("M112").Value = 2.25
$ws.Range("N112").Value = 3.3
$ws.Range("O112").Value = 2.7
$ws.Range("P112").Value = -0.25
$ws.Range("Q112").Value = 2.025
$ws.Range("R112").Value = 1.775
$ws.Range("S112").Value = 2.5
$ws.Range("T112").Value = 1.875
$ws.Range("U112").Value = 1.925
$ws.Range("V112").Value = 1.25
$ws.Range("W112").Value = -1
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = 1.025
$ws.Range("Z112").Value = -1
$ws.Range("AA112").Value = 0.875
$ws.Range("AB112").Value = -1

# Row 113
$ws.Range("A113").Value = 111
$ws.Range("B113").Value = 7802880
$ws.Range("C113").Value = "Canada Premier League"
$ws.Range("D113").Value = 45437.66666666666
$ws.Range("E113").Value = "Atletico Ottawa"
$ws.Range("F113").Value = "Forge FC"
$ws.Range("G113").Value = 3
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = "H"
$ws.Range("J113").Value = 2.375
$ws.Range("K113").Value = 3.2
$ws.Range("L113").Value = 2.625
$ws.Range("M113").Value = 2.4
$ws.Range("N113").Value = 3.2
$ws.Range("O113").Value = 2.6
$ws.Range("P113").Value = 0
$ws.Range("Q113").Value = 1.8
$ws.Range("R113").Value = 2
$ws.Range("S113").Value = 2.5
$ws.Range("T113").Value = 2
$ws.Range("U113").Value = 1.8
$ws.Range("V113").Value = 1.4
$ws.Range("W113").Value = -1
$ws.Range("X113").Value = -1
$ws.Range("Y113").Value = 0.8
$ws.Range("Z113").Value = -1
$ws.Range("AA113").Value = 1
$ws.Range("AB113").Value = -1

# Row 114
$ws.Range("A114").Value = 112
$ws.Range("B114").Value = 7802945
$ws.Range("C114").Value = "Canada Premier League"
$ws.Range("D114").Value = 45437.79166666666
$ws.Range("E114").Value = "Vancouver FC"
$ws.Range("F114").Value = "Pacific FC CA"
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 1
$ws.Range("I114").Value = "H"
$ws.Range("J114").Value = 2.75
$ws.Range("K114").Value = 3.25
$ws.Range("L114").Value = 2.25
$ws.Range("M114").Value = 2.55
$ws.Range("N114").Value = 3.5
$ws.Range("O114").Value = 2.375
$ws.Range("P114").Value = 0
$ws.Range("Q114").Value = 2
$ws.Range("R114").Value = 1.8
$ws.Range("S114").Value = 2.75
$ws.Range("T114").Value = 1.95
$ws.Range("U114").Value = 1.85
$ws.Range("V114").Value = 1.55
$ws.Range("W114").Value = -1
$ws.Range("X114").Value = -1
$ws.Range("Y114").Value = 1
$ws.Range("Z114").Value = -1
$ws.Range("AA114").Value = 0.475
$ws.Range("AB114").Value = -0.5

# Row 115
$ws.Range("A115").Value = 113
$ws.Range("B115").Value = 7803369
$ws.Range("C115").Value = "Canada Premier League"
$ws.Range("D115").Value = 45438.75
$ws.Range("E115").Value = "Cavalry FC"
$ws.Range("F115").Value = "Valour FC"
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 1
$ws.Range("I115").Value = "D"
$ws.Range("J115").Value = 1.571
$ws.Range("K115").Value = 3.6
$ws.Range("L115").Value = 5
$ws.Range("M115").Value = 1.42
$ws.Range("N115").Value = 4.2
$ws.Range("O115").Value = 6.5
$ws.Range("P115").Value = -1.25
$ws.Range("Q115").Value = 1.825
$ws.Range("R115").Value = 1.975
$ws.Range("S115").Value = 3
$ws.Range("T115").Value = 1.975
$ws.Range("U115").Value = 1.825
$ws.Range("V115").Value = -1
$ws.Range("W115").Value = 3.2
$ws.Range("X115").Value = -1
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.9750000000000001
$ws.Range("AA115").Value = -1
$ws.Range("AB115").Value = 0.825
